# Apply edit: replace the "Result" column text with a shorter summary,
# narrow column E to fit the new shorter text, and move the active
# selection from E8 to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Result" values (column E, rows 2-11) with the new shorter text.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "Outperformed market index "
}

# Column E is much narrower now that the text is shorter.
# (27.3 is used as input because this runtime's column-width rounding
# snaps to a 1/6-wide grid; it is the closest achievable value to the
# target stored width of 28.1640625.)
$ws.Range("E1:E11").ColumnWidth = 27.3

# Move the selection from E8 to G8.
$ws.Range("G8").Select()
